$d = $word.ActiveDocument

# Locate the "20 e 25" page-count figures inside the submission-norms paragraph
# ("O texto submetido deve conter entre 20 e 25 páginas, ...").
$whole = $d.Content
$found = $whole.Find.Execute("20 e 25", $true, $false, $false, $false, $false, `
                              $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find the '20 e 25' page-range text to update."
}
$start = $whole.Start
$end = $whole.End

# Update the text content first, while the sentence is still contained in a
# single run:
#   "20" -> "15"   (minimum page count)
$r20 = $d.Range($start, $start + 2)
$r20.Text = "15"

#   the "5" in "25" -> "0"   (maximum page count: 25 -> 20)
#   ("20" -> "15" keeps the same length, so $end is still valid here.)
$r5 = $d.Range($end - 1, $end)
$r5.Text = "0"

# Now split the edited fragments into their own runs, mirroring the
# run-boundary pattern Word leaves behind when text is retyped in place.
# Toggling Bold on/off forces the split without leaving any visible
# formatting change behind (the run's rPr, sz 24, ends up identical again).
$r20b = $d.Range($start, $start + 2)
$r20b.Bold = 1
$r20b.Bold = 0

$r5b = $d.Range($end - 1, $end)
$r5b.Bold = 1
$r5b.Bold = 0

Write-Output "Updated page-count range to 15 e 20."
